$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MatrixEdge")

$ws.Range("N1").Value = "num_references"
$ws.Range("O1").Value = "num_sentences"
